$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old content completely (old data ranged over A1:B56)
$ws.Cells.Clear()

# Data table to write: header row + 10 data rows, columns A:E
$data = @(
    @("mati", "fdsd", "dsfsf", "sfdfds", "sdffsd"),
    @(1, "fs", "sf", "sfd", "fs"),
    @(1, "dd", "fs", "fsdsdfs", "f"),
    @(1, "fs", "fs", "fsdsdfs", "fs"),
    @(2, "fs", "fs", "fsdsdfs", "fs"),
    @(3, "sf", "sf", "sfs", "f"),
    @(54, "sf", "sf", "fsdsdfs", "f"),
    @(5, "sf", "sf", "fsdsdfs", "s"),
    @(6, "sfd", "sf", "sfs", "fdfs"),
    @(4, "sf", "sf", "fsdsdfs", "fdfs"),
    @(5, "sdffsdf", "sf", "fsdsdfs", "fs")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Range("H12").Select()
